$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.280412554740906
$ws.Range("B1").Value = 1.202524185180664
$ws.Range("C1").Value = 1.028637528419495
$ws.Range("D1").Value = 1.086861729621887
$ws.Range("E1").Value = 1.005748271942139
